# This script reproduces the commit "Fruta / hortaliza, semanal" for the
# "Hortaliza, Terminal La Palmera de La Serena - Ajo" sheet.
#
# The edit inserts two new daily price rows (one new weekly observation,
# each unit of sale) right after the current first data block (row 60),
# pushing all the existing data down by two rows, and growing the used
# range from A1:R183 to A1:R185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 61-62; everything below (old rows 61..183)
# shifts down to rows 63..185.
$ws.Rows("61:62").Insert()

# --- Fill in the new row 61 ($/caja 10 kilos) ---
$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44536
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 100112003
$ws.Range("G61").Value = "Ajo"
$ws.Range("H61").Value = "Chino"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 19000
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = 19500
$ws.Range("N61").Value = "`$/caja 10 kilos"
$ws.Range("O61").Value = "China"
$ws.Range("P61").Value = 1950
$ws.Range("Q61").Value = 10
$ws.Range("R61").Value = "Hortaliza"

# --- Fill in the new row 62 ($/malla 10 kilos) ---
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44536
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112003
$ws.Range("G62").Value = "Ajo"
$ws.Range("H62").Value = "Chino"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = 20500
$ws.Range("N62").Value = "`$/malla 10 kilos"
$ws.Range("O62").Value = "China"
$ws.Range("P62").Value = 2050
$ws.Range("Q62").Value = 10
$ws.Range("R62").Value = "Hortaliza"
